$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet holds one weekly price record per row (Espinaca / Femacal de La
# Calera). The update adds two new weekly records:
#   - one inserted before the existing row 323 (pushes old 323.. down by 1)
#   - one inserted before what is now row 352, i.e. after the first insert
#     (pushes the remaining rows down by 1 more)
# All other columns (Mercado/Region/Categoria/etc.) are identical across the
# whole block, only Fecha (D), Volumen (J), Precio minimo/maximo/promedio
# (K/L/M) and Precio $/Kg (P) change per row.
# ---------------------------------------------------------------------------

function Set-EspinacaRow {
    param($Row, $Fecha, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg)

    $ws.Range("A$Row").Value = 3
    $ws.Range("B$Row").Value = "Femacal de La Calera"
    $ws.Range("C$Row").Value = "Coquimbo"
    $ws.Range("D$Row").Value = $Fecha
    $ws.Range("E$Row").Value = 5
    $ws.Range("F$Row").Value = 100112012
    $ws.Range("G$Row").Value = "Espinaca"
    $ws.Range("H$Row").Value = "Sin especificar"
    $ws.Range("I$Row").Value = "Primera"
    $ws.Range("J$Row").Value = $Volumen
    $ws.Range("K$Row").Value = $PrecioMin
    $ws.Range("L$Row").Value = $PrecioMax
    $ws.Range("M$Row").Value = $PrecioProm
    $ws.Range("N$Row").Value = "$/docena de atados (3 kilos)"
    $ws.Range("O$Row").Value = "Provincia de Quillota"
    $ws.Range("P$Row").Value = $PrecioKg
    $ws.Range("Q$Row").Value = 3
    $ws.Range("R$Row").Value = "Hortaliza"
}

# Insert new row before current row 323, shifting rows 323..362 down to 324..363.
$ws.Rows.Item(323).Insert()
Set-EspinacaRow 323 44748 230 4000 4500 4261 1420

# Insert another new row before (the now shifted) row 352, shifting rows
# 352..363 down to 353..364.
$ws.Rows.Item(352).Insert()
Set-EspinacaRow 352 44747 70 4000 4000 4000 1333
